$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 5
$ws.Range("H5").Value = 121.25
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()
# Row 100
$ws.Range("H100").Value = 1611
$ws.Range("I100").Value = 1071.4286
$ws.Range("K100").Value = 1071.4286
$ws.Range("M100").Value = -530.4286
# Row 103
$ws.Range("H103").Value = 499
$ws.Range("J103").Value = 499
$ws.Range("L103").Value = 1497
$ws.Range("N103").Value = -2669
# Row 112
$ws.Range("H112").Value = 2711.8
$ws.Range("J112").Value = 2902
$ws.Range("L112").Value = 8706
$ws.Range("N112").Value = -10922
# Row 132
$ws.Range("H132").Value = 1164.9584
$ws.Range("I132").Value = 1193.8695
$ws.Range("J132").Value = 500
$ws.Range("K132").Value = 3581.6085
$ws.Range("L132").Value = 1500
$ws.Range("M132").Value = -1051.6085
$ws.Range("N132").Value = -6560
# Row 138
$ws.Range("H138").Value = 3427.6567
$ws.Range("I138").Value = 2099.4285
$ws.Range("J138").Value = 3582.6167
$ws.Range("K138").Value = 6298.2855
$ws.Range("L138").Value = 10747.8501
$ws.Range("M138").Value = -1158.2855
$ws.Range("N138").Value = -21027.8501

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 10997.738
$ws.Range("I32").Value = 8664.306
$ws.Range("J32").Value = 24998.334
$ws.Range("K32").Value = 8664.306
$ws.Range("L32").Value = 24998.334
$ws.Range("M32").Value = -8377.306
$ws.Range("N32").Value = -25572.334
# Row 45
$ws.Range("H45").Value = 1721
$ws.Range("I45").Value = 1721
$ws.Range("K45").Value = 1721
$ws.Range("M45").Value = -1344
# Row 55
$ws.Range("H55").Value = 37332.332
$ws.Range("I55").Value = 2000
$ws.Range("J55").Value = 54998.5
$ws.Range("K55").Value = 2000
$ws.Range("L55").Value = 54998.5
$ws.Range("M55").Value = -1685
$ws.Range("N55").Value = -55628.5
# Row 61
$ws.Range("H61").Value = 5000
$ws.Range("I61").Value = 5000
$ws.Range("K61").Value = 5000
$ws.Range("M61").Value = -4788
# Row 88
$ws.Range("H88").Value = 187.1
$ws.Range("I88").Value = 191.22223
$ws.Range("K88").Value = 191.22223
$ws.Range("M88").Value = 214.77777
# Row 91
$ws.Range("H91").Value = 187.1
$ws.Range("I91").Value = 191.22223
$ws.Range("K91").Value = 191.22223
$ws.Range("M91").Value = 1212.77777
# Row 96
$ws.Range("H96").Value = 10000
$ws.Range("J96").Value = 10000
$ws.Range("L96").Value = 10000
$ws.Range("N96").Value = -15492
# Row 106
$ws.Range("H106").Value = 28500
$ws.Range("J106").Value = 28500
$ws.Range("L106").Value = 28500
$ws.Range("N106").Value = -31024
# Row 110
$ws.Range("H110").Value = 10532.667
$ws.Range("I110").Value = 10532.667
$ws.Range("K110").Value = 10532.667
$ws.Range("M110").Value = -8487.666999999999
# Row 132
$ws.Range("H132").Value = 3080
$ws.Range("I132").Value = 2730.9092
$ws.Range("K132").Value = 8192.7276
$ws.Range("M132").Value = -5662.7276
# Row 136
$ws.Range("H136").Value = 5000
$ws.Range("I136").Value = 5000
$ws.Range("K136").Value = 15000
$ws.Range("M136").Value = -12450
# Row 139
$ws.Range("H139").Value = 67857
$ws.Range("I139").Value = 54999
$ws.Range("J139").Value = 80715
$ws.Range("K139").Value = 54999
$ws.Range("L139").Value = 80715
$ws.Range("M139").Value = -49859
$ws.Range("N139").Value = -90995

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 33
$ws.Range("H33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()
# Row 86
$ws.Range("H86").Value = 992.375
$ws.Range("I86").Value = 947.7619
$ws.Range("J86").Value = 1304.6666
$ws.Range("K86").Value = 947.7619
$ws.Range("L86").Value = 1304.6666
$ws.Range("M86").Value = 175.2381
$ws.Range("N86").Value = -3550.6666
# Row 89
$ws.Range("H89").Value = 992.375
$ws.Range("I89").Value = 947.7619
$ws.Range("J89").Value = 1304.6666
$ws.Range("K89").Value = 4738.809499999999
$ws.Range("L89").Value = 6523.333000000001
$ws.Range("M89").Value = 877.1905000000006
$ws.Range("N89").Value = -17755.333
# Row 99
$ws.Range("H99").Value = 3162.1052
$ws.Range("I99").Value = 3193.1428
$ws.Range("J99").Value = 3075.2
$ws.Range("K99").Value = 3193.1428
$ws.Range("L99").Value = 3075.2
$ws.Range("M99").Value = -1695.1428
$ws.Range("N99").Value = -6071.2
# Row 107
$ws.Range("H107").Value = 724.2381
$ws.Range("I107").Value = 728.3333
$ws.Range("K107").Value = 728.3333
$ws.Range("M107").Value = 1191.6667

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 500
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 500
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 500
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -726
# Row 132
$ws.Range("H132").Value = 4795.5454
$ws.Range("I132").Value = 3541
$ws.Range("K132").Value = 10623
$ws.Range("M132").Value = -8093

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 2254280.5
$ws.Range("I4").Value = 3073046.2
$ws.Range("J4").Value = 2675
$ws.Range("K4").Value = 9219138.600000001
$ws.Range("L4").Value = 8025
$ws.Range("M4").Value = -9219026.600000001
$ws.Range("N4").Value = -8249
# Row 32
$ws.Range("H32").Value = 12305.556
$ws.Range("I32").Value = 1462.5
$ws.Range("K32").Value = 4387.5
$ws.Range("M32").Value = -4104.5
# Row 106
$ws.Range("H106").Value = 12285.429
$ws.Range("J106").Value = 12285.429
$ws.Range("L106").Value = 36856.287
$ws.Range("N106").Value = -38748.287

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 33
$ws.Range("H33").Value = 19482
$ws.Range("J33").Value = 19482
$ws.Range("L33").Value = 19482
$ws.Range("N33").Value = -19986
# Row 70
$ws.Range("H70").Value = 4900.4546
$ws.Range("I70").Value = 4726.25
$ws.Range("K70").Value = 4726.25
$ws.Range("M70").Value = -4456.25
# Row 73
$ws.Range("H73").Value = 4900.4546
$ws.Range("I73").Value = 4726.25
$ws.Range("K73").Value = 4726.25
$ws.Range("M73").Value = -3790.25
# Row 97
$ws.Range("H97").Value = 753.6667
$ws.Range("I97").Value = 620.6875
$ws.Range("K97").Value = 620.6875
$ws.Range("M97").Value = -124.6875
# Row 109
$ws.Range("H109").Value = 41944.445
$ws.Range("J109").Value = 41944.445
$ws.Range("L109").Value = 41944.445
$ws.Range("N109").Value = -44024.445
# Row 132
$ws.Range("H132").Value = 1974.909
$ws.Range("I132").Value = 1193.3158
$ws.Range("K132").Value = 3579.9474
$ws.Range("M132").Value = -1049.9474

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 1071.6364
$ws.Range("J46").Value = 1337.5
$ws.Range("L46").Value = 1337.5
$ws.Range("N46").Value = -1713.5
# Row 64
$ws.Range("H64").Value = 32500
$ws.Range("J64").Value = 32500
$ws.Range("L64").Value = 32500
$ws.Range("N64").Value = -32950
# Row 67
$ws.Range("H67").Value = 32500
$ws.Range("J67").Value = 32500
$ws.Range("L67").Value = 32500
$ws.Range("N67").Value = -34060
# Row 122
$ws.Range("H122").Value = 4166.1816
$ws.Range("I122").Value = 4182.9
$ws.Range("K122").Value = 12548.7
$ws.Range("M122").Value = -10098.7
# Row 132
$ws.Range("H132").Value = 7000
$ws.Range("J132").Value = 10000
$ws.Range("L132").Value = 30000
$ws.Range("N132").Value = -35060
# Row 133
$ws.Range("H133").Value = 66333.336
$ws.Range("J133").Value = 66333.336
$ws.Range("L133").Value = 66333.336
$ws.Range("N133").Value = -71393.336

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 2977.5
$ws.Range("I81").Value = 2527.6667
$ws.Range("K81").Value = 5055.3334
$ws.Range("M81").Value = -3994.3334
# Row 84
$ws.Range("H84").Value = 2977.5
$ws.Range("I84").Value = 2527.6667
$ws.Range("K84").Value = 25276.667
$ws.Range("M84").Value = -19972.667
# Row 104
$ws.Range("H104").Value = 26749.75
$ws.Range("J104").Value = 26749.75
$ws.Range("L104").Value = 26749.75
$ws.Range("N104").Value = -33737.75
# Row 107
$ws.Range("H107").Value = 619.5
$ws.Range("I107").Value = 436.5
$ws.Range("K107").Value = 1309.5
$ws.Range("M107").Value = 610.5
# Row 113
$ws.Range("H113").Value = 717.75
$ws.Range("J113").Value = 297.75
$ws.Range("L113").Value = 893.25
$ws.Range("N113").Value = -5233.25
# Row 132
$ws.Range("H132").Value = 1358.875
$ws.Range("I132").Value = 1231.5714
$ws.Range("K132").Value = 3694.7142
$ws.Range("M132").Value = -1164.7142
# Row 136
$ws.Range("H136").Value = 5595.3
$ws.Range("I136").Value = 5244.75
$ws.Range("K136").Value = 15734.25
$ws.Range("M136").Value = -13184.25
